$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Ensure the data in columns D and E stays formatted as text, matching the
# original inline-string cell values (e.g. "61.914.88" must not turn into a number/date).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "62.400.84"
$ws.Cells.Item(2, 5).Value = "  -1.02%  "
$ws.Cells.Item(3, 4).Value = "3.029.24"
$ws.Cells.Item(3, 5).Value = "  -1.30%  "
$ws.Cells.Item(4, 4).Value = "0.999"
$ws.Cells.Item(4, 5).Value = "  -0.12%  "
$ws.Cells.Item(5, 4).Value = "545.34"
$ws.Cells.Item(5, 5).Value = "  +1.25%  "
$ws.Cells.Item(6, 4).Value = "135.04"
$ws.Cells.Item(6, 5).Value = "  -1.54%  "
$ws.Cells.Item(7, 5).Value = "  +0.02%  "
$ws.Cells.Item(8, 4).Value = "3.022.93"
$ws.Cells.Item(8, 5).Value = "  -1.22%  "
$ws.Cells.Item(9, 4).Value = "0.494"
$ws.Cells.Item(9, 5).Value = "  +0.35%  "
$ws.Cells.Item(10, 4).Value = "6.34"
$ws.Cells.Item(10, 5).Value = "  +1.94%  "
$ws.Cells.Item(11, 4).Value = "0.148"
$ws.Cells.Item(11, 5).Value = "  -4.69%  "
$ws.Cells.Item(12, 4).Value = "0.450"
$ws.Cells.Item(12, 5).Value = "  -0.72%  "
$ws.Cells.Item(13, 4).Value = "34.94"
$ws.Cells.Item(13, 5).Value = "  +1.58%  "
$ws.Cells.Item(14, 4).Value = "0.0000223"
$ws.Cells.Item(14, 5).Value = "  +0.56%  "
$ws.Cells.Item(15, 4).Value = "3.504.03"
$ws.Cells.Item(15, 5).Value = "  -1.54%  "
$ws.Cells.Item(16, 4).Value = "62.201.14"
$ws.Cells.Item(16, 5).Value = "  -1.29%  "
$ws.Cells.Item(17, 5).Value = "  -2.59%  "
$ws.Cells.Item(18, 4).Value = "3.017.57"
$ws.Cells.Item(18, 5).Value = "  -1.51%  "
$ws.Cells.Item(19, 4).Value = "6.70"
$ws.Cells.Item(19, 5).Value = "  +1.11%  "
$ws.Cells.Item(20, 4).Value = "484.90"
$ws.Cells.Item(20, 5).Value = "  +3.39%  "
$ws.Cells.Item(21, 4).Value = "13.35"
$ws.Cells.Item(21, 5).Value = "  -1.04%  "
$ws.Cells.Item(22, 4).Value = "0.676"
$ws.Cells.Item(22, 5).Value = "  -2.58%  "
$ws.Cells.Item(23, 4).Value = "7.05"
$ws.Cells.Item(23, 5).Value = "  +0.54%  "
$ws.Cells.Item(24, 4).Value = "82.11"
$ws.Cells.Item(24, 5).Value = "  +4.73%  "
$ws.Cells.Item(25, 4).Value = "12.13"
$ws.Cells.Item(25, 5).Value = "  +0.21%  "
$ws.Cells.Item(26, 5).Value = "  +0.07%  "
$ws.Cells.Item(27, 4).Value = "2.71"
$ws.Cells.Item(27, 5).Value = "  +0.83%  "
$ws.Cells.Item(28, 4).Value = "7.83"
$ws.Cells.Item(28, 5).Value = "  -0.35%  "
$ws.Cells.Item(29, 4).Value = "0.999"
$ws.Cells.Item(29, 5).Value = "  -0.08%  "
$ws.Cells.Item(30, 4).Value = "1.94"
$ws.Cells.Item(30, 5).Value = "  +3.46%  "
$ws.Cells.Item(31, 4).Value = "25.95"
$ws.Cells.Item(31, 5).Value = "  -0.50%  "
$ws.Cells.Item(32, 4).Value = "1.14"
$ws.Cells.Item(32, 5).Value = "  -1.58%  "
$ws.Cells.Item(33, 4).Value = "5.69"
$ws.Cells.Item(33, 5).Value = "  +4.59%  "
$ws.Cells.Item(34, 4).Value = "2.38"
$ws.Cells.Item(34, 5).Value = "  +3.25%  "
$ws.Cells.Item(35, 4).Value = "55.57"
$ws.Cells.Item(35, 5).Value = "  -5.64%  "
$ws.Cells.Item(36, 4).Value = "5.92"
$ws.Cells.Item(36, 5).Value = "  -0.47%  "
$ws.Cells.Item(37, 4).Value = "452.89"
$ws.Cells.Item(37, 5).Value = "  -5.70%  "
$ws.Cells.Item(38, 4).Value = "3.167.32"
$ws.Cells.Item(38, 5).Value = "  -2.63%  "
$ws.Cells.Item(39, 4).Value = "0.0803"
$ws.Cells.Item(39, 5).Value = "  +1.50%  "
$ws.Cells.Item(40, 4).Value = "0.0389"
$ws.Cells.Item(40, 5).Value = "  -2.02%  "
$ws.Cells.Item(41, 5).Value = "  +1.16%  "
$ws.Cells.Item(42, 4).Value = "8.16"
$ws.Cells.Item(42, 5).Value = "  +0.50%  "
$ws.Cells.Item(43, 4).Value = "2.47"
$ws.Cells.Item(43, 5).Value = "  -3.69%  "
$ws.Cells.Item(44, 4).Value = "26.45"
$ws.Cells.Item(44, 5).Value = "  +5.03%  "
$ws.Cells.Item(45, 5).Value = "  -0.08%  "
$ws.Cells.Item(46, 4).Value = "0.246"
$ws.Cells.Item(46, 5).Value = "  -1.85%  "
$ws.Cells.Item(47, 4).Value = "0.110"
$ws.Cells.Item(47, 5).Value = "  +1.08%  "
$ws.Cells.Item(48, 4).Value = "1.99"
$ws.Cells.Item(48, 5).Value = "  -0.71%  "
$ws.Cells.Item(49, 4).Value = "116.59"
$ws.Cells.Item(49, 5).Value = "  -5.63%  "
$ws.Cells.Item(50, 5).Value = "  +4.43%  "
$ws.Cells.Item(51, 4).Value = "0.0₃0495"
$ws.Cells.Item(51, 5).Value = "  -4.59%  "
